$d = $word.ActiveDocument

$replacements = @(
    @{old = "397÷4=99, 1";  new = "731÷5=146, 1"},
    @{old = "911÷2=455, 1"; new = "301÷7=43, 0"},
    @{old = "150÷9=16, 6";  new = "487÷2=243, 1"},
    @{old = "826÷2=413, 0"; new = "778÷7=111, 1"},
    @{old = "807÷9=89, 6";  new = "948÷6=158, 0"},
    @{old = "491÷5=98, 1";  new = "289÷4=72, 1"},
    @{old = "744÷9=82, 6";  new = "490÷8=61, 2"},
    @{old = "130÷2=65, 0";  new = "564÷4=141, 0"},
    @{old = "782÷2=391, 0"; new = "550÷7=78, 4"},
    @{old = "539÷9=59, 8";  new = "164÷9=18, 2"},
    @{old = "552÷2=276, 0"; new = "872÷2=436, 0"},
    @{old = "213÷5=42, 3";  new = "151÷9=16, 7"},
    @{old = "270÷2=135, 0"; new = "165÷6=27, 3"},
    @{old = "400÷3=133, 1"; new = "225÷3=75, 0"},
    @{old = "913÷6=152, 1"; new = "508÷5=101, 3"},
    @{old = "415÷2=207, 1"; new = "538÷3=179, 1"},
    @{old = "821÷9=91, 2";  new = "668÷4=167, 0"},
    @{old = "963÷4=240, 3"; new = "507÷6=84, 3"},
    @{old = "426÷9=47, 3";  new = "674÷5=134, 4"},
    @{old = "763÷5=152, 3"; new = "222÷7=31, 5"},
    @{old = "246÷5=49, 1";  new = "491÷4=122, 3"},
    @{old = "661÷9=73, 4";  new = "330÷4=82, 2"},
    @{old = "796÷4=199, 0"; new = "722÷9=80, 2"},
    @{old = "991÷2=495, 1"; new = "157÷9=17, 4"},
    @{old = "928÷3=309, 1"; new = "992÷4=248, 0"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null
}
